# Select the 11 flow-chart / arrow shapes on slide 4 (the benchmarking
# diagram: "Read Crawled Documents" -> ... -> "Save Model") and group them
# into a single group shape, matching PowerPoint's native "Group" command.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$range = $s.Shapes.Range(@(1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11))
$grp = $range.Group()
